$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F2:F21").Value = "s"
